# Update odds values on Sheet1 as scraped from a later FlashScore run.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "G5"  = 1.65
    "G6"  = 2.45
    "R6"  = 1.53
    "G11" = 3.8
    "H11" = 3.1
    "I11" = 1.98
    "J11" = 4.3
    "L11" = 2.62
    "R11" = 1.53
    "S11" = 1.47
    "T11" = 2.35
    "U11" = 1.98
    "V11" = 1.65
    "W11" = 9
    "X11" = 19
    "Y11" = 13.5
    "Z11" = 60
    "AC11" = 7.4
    "AG11" = 5.8
    "AH11" = 8.5
    "AJ11" = 17.5
    "AN11" = 5.4
    "AP11" = 32
    "AU11" = 7.6
    "AW11" = 3.7
    "AX11" = 10.25
    "AY11" = 22
    "BA11" = 90
    "BB11" = 350
    "AC12" = 11.75
    "AP12" = 16
    "AT12" = 2.9
    "AX12" = 23
    "AY12" = 27
    "AE13" = 11
    "AH13" = 19
    "AJ13" = 34
    "AN13" = 4.5
    "AO13" = 11
    "AR13" = 41
    "BC13" = 301
    "K19" = 1.91
    "G21" = 1.7
    "J21" = 2.38
    "Q21" = 2.08
    "R21" = 1.73
    "AC21" = 9
    "AD21" = 7
    "AM21" = 351
    "AU21" = 8.5
    "AX21" = 29
    "R22" = 1.67
    "I23" = 2.88
    "R23" = 1.65
    "R24" = 1.75
    "R25" = 1.65
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
